$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 430: remove now-empty placeholder cells B430:E430 ---
$ws.Range("B430:E430").ClearContents()

# --- Row 431 ---
$ws.Cells.Item(431, 1).Value = 'The Temple of Light'
$c = $ws.Cells.Item(431, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(431, 7)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$ws.Cells.Item(431, 8).Value = 60
$ws.Cells.Item(431, 9).Value = 'Vicki'
$ws.Cells.Item(431, 13).Value = 'Jonathan Morris'
$ws.Cells.Item(431, 14).Value = 'Lisa Bowerman'
$c = $ws.Cells.Item(431, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 432 ---
$ws.Cells.Item(432, 1).Value = 'Stardust and Ashes'
$c = $ws.Cells.Item(432, 6)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$c = $ws.Cells.Item(432, 7)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$ws.Cells.Item(432, 8).Value = 76
$ws.Cells.Item(432, 9).Value = 'Susan'
$ws.Cells.Item(432, 11).Value = '1st Dr, Barbara, Ian'
$ws.Cells.Item(432, 13).Value = 'Ian Potter'
$ws.Cells.Item(432, 14).Value = 'Lisa Bowerman'
$c = $ws.Cells.Item(432, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 433 ---
$ws.Cells.Item(433, 1).Value = 'The White Ship'
$c = $ws.Cells.Item(433, 6)
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"
$c = $ws.Cells.Item(433, 7)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$ws.Cells.Item(433, 8).Value = 78
$ws.Cells.Item(433, 9).Value = '1st Dr'
$ws.Cells.Item(433, 10).Value = 'Steven'
$ws.Cells.Item(433, 13).Value = 'Paul Morris'
$ws.Cells.Item(433, 14).Value = 'Lisa Bowerman'
$c = $ws.Cells.Item(433, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 434 ---
$ws.Cells.Item(434, 1).Value = 'The Y Factor'
$c = $ws.Cells.Item(434, 6)
$c.NumberFormat = "@"
$c.Value = '4'
$c.Style = "Normal"
$c = $ws.Cells.Item(434, 7)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$ws.Cells.Item(434, 8).Value = 71
$ws.Cells.Item(434, 9).Value = '1st Dr'
$ws.Cells.Item(434, 10).Value = 'Dodo'
$ws.Cells.Item(434, 12).Value = 'Sanderson, a fungus'
$ws.Cells.Item(434, 13).Value = 'Christopher Cooper'
$ws.Cells.Item(434, 14).Value = 'Lisa Bowerman'
$c = $ws.Cells.Item(434, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 435 ---
$ws.Cells.Item(435, 1).Value = 'Dark Watchers of California'
$c = $ws.Cells.Item(435, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(435, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(435, 8).Value = 36
$ws.Cells.Item(435, 9).Value = '12th Dr'
$ws.Cells.Item(435, 10).Value = 'Bill Potts'
$ws.Cells.Item(435, 12).Value = 'Finch'
$ws.Cells.Item(435, 13).Value = 'Riley Silverman'
$ws.Cells.Item(435, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(435, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 436 ---
$ws.Cells.Item(436, 1).Value = 'Dark Watchers of California'
$c = $ws.Cells.Item(436, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(436, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(436, 8).Value = 36
$ws.Cells.Item(436, 9).Value = '12th Dr'
$ws.Cells.Item(436, 10).Value = 'Bill Potts'
$ws.Cells.Item(436, 12).Value = 'Finch'
$ws.Cells.Item(436, 13).Value = 'Riley Silverman'
$ws.Cells.Item(436, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(436, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 437 ---
$ws.Cells.Item(437, 1).Value = 'Dark Watchers of California'
$c = $ws.Cells.Item(437, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(437, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(437, 8).Value = 36
$ws.Cells.Item(437, 9).Value = '12th Dr'
$ws.Cells.Item(437, 10).Value = 'Bill Potts'
$ws.Cells.Item(437, 12).Value = 'Finch'
$ws.Cells.Item(437, 13).Value = 'Riley Silverman'
$ws.Cells.Item(437, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(437, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 438 ---
$ws.Cells.Item(438, 1).Value = 'Dark Watchers of California'
$c = $ws.Cells.Item(438, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(438, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(438, 8).Value = 36
$ws.Cells.Item(438, 9).Value = '12th Dr'
$ws.Cells.Item(438, 10).Value = 'Bill Potts'
$ws.Cells.Item(438, 12).Value = 'Finch'
$ws.Cells.Item(438, 13).Value = 'Riley Silverman'
$ws.Cells.Item(438, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(438, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 439 ---
$ws.Cells.Item(439, 1).Value = 'When I Say Run'
$c = $ws.Cells.Item(439, 6)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$ws.Cells.Item(439, 7).Value = 'N/A'
$ws.Cells.Item(439, 8).Value = 41
$ws.Cells.Item(439, 9).Value = 'N, /, A'
$ws.Cells.Item(439, 10).Value = 'N/A'
$ws.Cells.Item(439, 11).Value = 'N, /, A'
$ws.Cells.Item(439, 12).Value = 'N/A'
$ws.Cells.Item(439, 13).Value = 'N/A'
$ws.Cells.Item(439, 14).Value = 'N/A'
$c = $ws.Cells.Item(439, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 440 ---
$ws.Cells.Item(440, 1).Value = 'Rise of the Eukaryans'
$c = $ws.Cells.Item(440, 6)
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"
$c = $ws.Cells.Item(440, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(440, 8).Value = 38
$ws.Cells.Item(440, 9).Value = '11th Dr'
$ws.Cells.Item(440, 12).Value = 'Eukaryans'
$ws.Cells.Item(440, 13).Value = 'Daniel Hardcastle'
$ws.Cells.Item(440, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(440, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 441 ---
$ws.Cells.Item(441, 1).Value = 'Ahead of Time'
$c = $ws.Cells.Item(441, 6)
$c.NumberFormat = "@"
$c.Value = '4'
$c.Style = "Normal"
$ws.Cells.Item(441, 7).Value = 'N/A'
$ws.Cells.Item(441, 8).Value = 45
$ws.Cells.Item(441, 9).Value = 'N, /, A'
$ws.Cells.Item(441, 10).Value = 'N/A'
$ws.Cells.Item(441, 11).Value = 'N, /, A'
$ws.Cells.Item(441, 12).Value = 'N/A'
$ws.Cells.Item(441, 13).Value = 'N/A'
$ws.Cells.Item(441, 14).Value = 'N/A'
$c = $ws.Cells.Item(441, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 442 ---
$ws.Cells.Item(442, 1).Value = 'Emerald Isle'
$c = $ws.Cells.Item(442, 6)
$c.NumberFormat = "@"
$c.Value = '5'
$c.Style = "Normal"
$ws.Cells.Item(442, 7).Value = 'N/A'
$ws.Cells.Item(442, 8).Value = 49
$ws.Cells.Item(442, 9).Value = 'N, /, A'
$ws.Cells.Item(442, 10).Value = 'N/A'
$ws.Cells.Item(442, 11).Value = 'N, /, A'
$ws.Cells.Item(442, 12).Value = 'N/A'
$ws.Cells.Item(442, 13).Value = 'N/A'
$ws.Cells.Item(442, 14).Value = 'N/A'
$c = $ws.Cells.Item(442, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 443 ---
$ws.Cells.Item(443, 1).Value = 'Dark is the Devil that Walks'
$c = $ws.Cells.Item(443, 6)
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$ws.Cells.Item(443, 7).Value = 'N/A'
$ws.Cells.Item(443, 8).Value = 44
$ws.Cells.Item(443, 9).Value = 'N, /, A'
$ws.Cells.Item(443, 10).Value = 'N/A'
$ws.Cells.Item(443, 11).Value = 'N, /, A'
$ws.Cells.Item(443, 12).Value = 'N/A'
$ws.Cells.Item(443, 13).Value = 'N/A'
$ws.Cells.Item(443, 14).Value = 'N/A'
$c = $ws.Cells.Item(443, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 444 ---
$ws.Cells.Item(444, 1).Value = 'Hooklight 1'
$c = $ws.Cells.Item(444, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(444, 7)
$c.NumberFormat = "@"
$c.Value = '12'
$c.Style = "Normal"
$ws.Cells.Item(444, 8).Value = 216
$ws.Cells.Item(444, 9).Value = '5th Dr'
$ws.Cells.Item(444, 10).Value = 'Adric, Nyssa, Tegan'
$ws.Cells.Item(444, 11).Value = '8th Dr'
$ws.Cells.Item(444, 12).Value = 'Nura, Nigh Guard'
$ws.Cells.Item(444, 13).Value = 'Tim Foley'
$ws.Cells.Item(444, 14).Value = 'Ken Bentley'
$c = $ws.Cells.Item(444, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 445 ---
$ws.Cells.Item(445, 1).Value = 'Missy Part 1'
$c = $ws.Cells.Item(445, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(445, 7).Value = 'N/A'
$ws.Cells.Item(445, 8).Value = 78
$ws.Cells.Item(445, 9).Value = 'N, /, A'
$ws.Cells.Item(445, 10).Value = 'N/A'
$ws.Cells.Item(445, 11).Value = 'N, /, A'
$ws.Cells.Item(445, 12).Value = 'N/A'
$ws.Cells.Item(445, 13).Value = 'N/A'
$ws.Cells.Item(445, 14).Value = 'N/A'
$c = $ws.Cells.Item(445, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 446 ---
$ws.Cells.Item(446, 1).Value = 'Dark Watchers of California'
$c = $ws.Cells.Item(446, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(446, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(446, 8).Value = 36
$ws.Cells.Item(446, 9).Value = '12th Dr'
$ws.Cells.Item(446, 10).Value = 'Bill Potts'
$ws.Cells.Item(446, 12).Value = 'Finch'
$ws.Cells.Item(446, 13).Value = 'Riley Silverman'
$ws.Cells.Item(446, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(446, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 447 ---
$ws.Cells.Item(447, 1).Value = 'When I Say Run'
$c = $ws.Cells.Item(447, 6)
$c.NumberFormat = "@"
$c.Value = '2'
$c.Style = "Normal"
$ws.Cells.Item(447, 7).Value = 'N/A'
$ws.Cells.Item(447, 8).Value = 41
$ws.Cells.Item(447, 9).Value = 'N, /, A'
$ws.Cells.Item(447, 10).Value = 'N/A'
$ws.Cells.Item(447, 11).Value = 'N, /, A'
$ws.Cells.Item(447, 12).Value = 'N/A'
$ws.Cells.Item(447, 13).Value = 'N/A'
$ws.Cells.Item(447, 14).Value = 'N/A'
$c = $ws.Cells.Item(447, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 448 ---
$ws.Cells.Item(448, 1).Value = 'Rise of the Eukaryans'
$c = $ws.Cells.Item(448, 6)
$c.NumberFormat = "@"
$c.Value = '3'
$c.Style = "Normal"
$c = $ws.Cells.Item(448, 7)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(448, 8).Value = 38
$ws.Cells.Item(448, 9).Value = '11th Dr'
$ws.Cells.Item(448, 12).Value = 'Eukaryans'
$ws.Cells.Item(448, 13).Value = 'Daniel Hardcastle'
$ws.Cells.Item(448, 14).Value = 'Peter Anghelides'
$c = $ws.Cells.Item(448, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 449 ---
$ws.Cells.Item(449, 1).Value = 'Ahead of Time'
$c = $ws.Cells.Item(449, 6)
$c.NumberFormat = "@"
$c.Value = '4'
$c.Style = "Normal"
$ws.Cells.Item(449, 7).Value = 'N/A'
$ws.Cells.Item(449, 8).Value = 45
$ws.Cells.Item(449, 9).Value = 'N, /, A'
$ws.Cells.Item(449, 10).Value = 'N/A'
$ws.Cells.Item(449, 11).Value = 'N, /, A'
$ws.Cells.Item(449, 12).Value = 'N/A'
$ws.Cells.Item(449, 13).Value = 'N/A'
$ws.Cells.Item(449, 14).Value = 'N/A'
$c = $ws.Cells.Item(449, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 450 ---
$ws.Cells.Item(450, 1).Value = 'Emerald Isle'
$c = $ws.Cells.Item(450, 6)
$c.NumberFormat = "@"
$c.Value = '5'
$c.Style = "Normal"
$ws.Cells.Item(450, 7).Value = 'N/A'
$ws.Cells.Item(450, 8).Value = 49
$ws.Cells.Item(450, 9).Value = 'N, /, A'
$ws.Cells.Item(450, 10).Value = 'N/A'
$ws.Cells.Item(450, 11).Value = 'N, /, A'
$ws.Cells.Item(450, 12).Value = 'N/A'
$ws.Cells.Item(450, 13).Value = 'N/A'
$ws.Cells.Item(450, 14).Value = 'N/A'
$c = $ws.Cells.Item(450, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 451 ---
$ws.Cells.Item(451, 1).Value = 'Dark is the Devil that Walks'
$c = $ws.Cells.Item(451, 6)
$c.NumberFormat = "@"
$c.Value = '6'
$c.Style = "Normal"
$ws.Cells.Item(451, 7).Value = 'N/A'
$ws.Cells.Item(451, 8).Value = 44
$ws.Cells.Item(451, 9).Value = 'N, /, A'
$ws.Cells.Item(451, 10).Value = 'N/A'
$ws.Cells.Item(451, 11).Value = 'N, /, A'
$ws.Cells.Item(451, 12).Value = 'N/A'
$ws.Cells.Item(451, 13).Value = 'N/A'
$ws.Cells.Item(451, 14).Value = 'N/A'
$c = $ws.Cells.Item(451, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 452 ---
$ws.Cells.Item(452, 1).Value = 'Hooklight 1'
$c = $ws.Cells.Item(452, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$c = $ws.Cells.Item(452, 7)
$c.NumberFormat = "@"
$c.Value = '12'
$c.Style = "Normal"
$ws.Cells.Item(452, 8).Value = 216
$ws.Cells.Item(452, 9).Value = '5th Dr'
$ws.Cells.Item(452, 10).Value = 'Adric, Nyssa, Tegan'
$ws.Cells.Item(452, 11).Value = '8th Dr'
$ws.Cells.Item(452, 12).Value = 'Nura, Nigh Guard'
$ws.Cells.Item(452, 13).Value = 'Tim Foley'
$ws.Cells.Item(452, 14).Value = 'Ken Bentley'
$c = $ws.Cells.Item(452, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

# --- Row 453 ---
$ws.Cells.Item(453, 1).Value = 'Missy Part 1'
$c = $ws.Cells.Item(453, 6)
$c.NumberFormat = "@"
$c.Value = '1'
$c.Style = "Normal"
$ws.Cells.Item(453, 7).Value = 'N/A'
$ws.Cells.Item(453, 8).Value = 78
$ws.Cells.Item(453, 9).Value = 'N, /, A'
$ws.Cells.Item(453, 10).Value = 'N/A'
$ws.Cells.Item(453, 11).Value = 'N, /, A'
$ws.Cells.Item(453, 12).Value = 'N/A'
$ws.Cells.Item(453, 13).Value = 'N/A'
$ws.Cells.Item(453, 14).Value = 'N/A'
$c = $ws.Cells.Item(453, 15)
$c.NumberFormat = "@"
$c.Value = '2025'
$c.Style = "Normal"

